# Replaces the "thermochemical water splitting" hydrogen production pathway
# with "hydrocarbon partial oxidation" on the HPtFM sheet, and removes the
# now-obsolete explanatory notes about thermochemical water splitting on the
# About sheet.

$wb = $excel.ActiveWorkbook

# --- About sheet ---
$wsAbout = $wb.Worksheets.Item("About")

# Remove the three note rows (14-16) that explained thermochemical water
# splitting's lack of tracked fuel use; they no longer apply.
$wsAbout.Range("A14").EntireRow.Delete() | Out-Null
$wsAbout.Range("A14").EntireRow.Delete() | Out-Null
$wsAbout.Range("A14").EntireRow.Delete() | Out-Null

# --- HPtFM sheet ---
$wsHPtFM = $wb.Worksheets.Item("HPtFM")

# Rename the pathway from thermochemical water splitting to hydrocarbon
# partial oxidation.
$wsHPtFM.Range("A6").Value = "hydrocarbon partial oxidation"

# Hydrocarbon partial oxidation consumes heavy/residual fuel oil (column I),
# whereas thermochemical water splitting used no tracked fuel at all.
$wsHPtFM.Range("I6").Value = 1

$wsHPtFM.Activate() | Out-Null
